# Add two new columns (I: "I0", J: "IF") to Sheet1, mirroring the header
# style already used by the existing columns (A1:H1) and filling in the
# per-row numeric values for rows 2-20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (row 1) -- text labels, matching the bold/bordered header
# style already applied to A1:H1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting from H1 onto the two new header cells so they
# pick up the same style (s="1") as the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Per-row data for the new "I0" (col I) and "IF" (col J) columns.
$data = @{
    2  = @(1, 5)
    3  = @(7, 7)
    4  = @(8, 8)
    5  = @(3, 5)
    6  = @(9, 9)
    7  = @(1, 2)
    8  = @(7, 7)
    9  = @(6, 6)
    10 = @(8, 8)
    11 = @(6, 7)
    12 = @(1, 7)
    13 = @(1, 5)
    14 = @(1, 6)
    15 = @(1, 5)
    16 = @(5, 5)
    17 = @(7, 7)
    18 = @(3, 4)
    19 = @(1, 3)
    20 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 9).Value = $values[0]
    $ws.Cells.Item($row, 10).Value = $values[1]
}
